$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")
$ws.Select()
$ws.Range("M1").Value = "CLIENT_GROUPS"
$ws.Range("A2").Select()
